$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells stay text (not auto-converted to numbers/dates) by forcing Text format first
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.71%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "45.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.66%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.698"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.41%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08365"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.22%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.91%"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9831"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.64%"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.597"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.17%"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1163"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.34%"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1945"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.29%"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "10.39"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-17.38%"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1011"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.92%"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04668"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.18%"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.64%"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001297"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.78%"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006057"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.18%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.368"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.05%"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.474"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.39%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3351"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-3.67%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.41%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2594"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.98%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04197"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.69%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001310"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.23%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004584"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "6.19%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.66%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003742"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.05%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02784"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "8.83%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05822"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.55%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007726"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.90%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.80%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007202"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-5.34%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001973"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.99%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008170"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.02%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007210"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.49%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.16%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005804"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.13%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003496"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "89.93%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003500"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.75%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.16%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.16%"
